$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
  2 = @("ECs", "Wnt5a", "Ryk", "ECs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 11.49304366666667, 34.479131, 0.2201107113282896, 0.2201107113282896, 1.837741513314555, 16.539673619831, 0.0058280296591177, 0.005828029659117701)
  3 = @("ECs", "Wnt5a", "Ryk", "FAPs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 27.87722133333333, 83.631664, 0.5338946927812221, 0.5338946927812221, 4.457576983607112, 40.118192852464, 0.01413631388312443, 0.01413631388312444)
  4 = @("ECs", "Wnt5a", "Ryk", "M1", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 0.133321, 0.399963, 0.002553316684083387, 0.002553316684083387, 0.02131807234033333, 0.191862651063, 0.0000676060027890405, 0.00006760600278904051)
  5 = @("ECs", "Wnt5a", "Ryk", "M2", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 2, 0.6666666666666666, 0.1029413333333333, 0.308824, 0.00197149604249735, 0.00197149604249735, 0.01646035351377778, 0.148143181624, 0.00005220071907982149, 0.00005220071907982149)
  6 = @("ECs", "Wnt5a", "Ryk", "sCs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 12.608304, 37.824912, 0.2414697831639074, 0.2414697831639074, 2.016072012368, 18.144648111312, 0.006393569170566306, 0.006393569170566307)
  7 = @("FAPs", "Wnt5a", "Ryk", "ECs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 11.49304366666667, 34.479131, 0.2201107113282896, 0.2201107113282896, 67.56935065897066, 608.1241559307359, 0.2142826816691719, 0.2142826816691719)
  8 = @("FAPs", "Wnt5a", "Ryk", "FAPs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 27.87722133333333, 83.631664, 0.5338946927812221, 0.5338946927812221, 163.8944215563093, 1475.049794006784, 0.5197583788980977, 0.5197583788980977)
  9 = @("FAPs", "Wnt5a", "Ryk", "M1", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 0.133321, 0.399963, 0.002553316684083387, 0.002553316684083387, 0.7838144237920001, 7.054329814128001, 0.002485710681294346, 0.002485710681294347)
  10 = @("FAPs", "Wnt5a", "Ryk", "M2", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 2, 0.6666666666666666, 0.1029413333333333, 0.308824, 0.00197149604249735, 0.00197149604249735, 0.6052077457493333, 5.446869711744, 0.001919295323417529, 0.001919295323417529)
  11 = @("FAPs", "Wnt5a", "Ryk", "sCs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 12.608304, 37.824912, 0.2414697831639074, 0.2414697831639074, 74.126135678208, 667.1352211038719, 0.2350762139933411, 0.2350762139933411)
}

foreach ($r in $rows.Keys) {
  $vals = $rows[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
  }
}